$d = $word.ActiveDocument

# Locate the paragraph that holds the email address "gam4kv@umsystem.edu"
$findRng = $d.Content
$found = $findRng.Find.Execute("gam4kv@umsystem.edu", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the email address paragraph"
}
$emailParaIndex = $findRng.Paragraphs.First.Index

# Re-fetch the paragraph by index so we get the full paragraph range
# (including its end-of-paragraph mark) rather than the collapsed Find range.
$emailPara = $d.Paragraphs.Item($emailParaIndex)
$emailRange = $emailPara.Range

# Make the paragraph mark match the existing (bold, 9pt) run formatting -
# this is what Word does automatically when you press Enter at the end of
# such a paragraph, recording the mark's rPr inside pPr.
$emailRange.Font.Bold = 1
$emailRange.Font.Size = 9

# Insert a brand-new paragraph right after the email paragraph; it inherits
# the same paragraph/run formatting.
$emailRange.InsertParagraphAfter()

# The newly-created paragraph is now immediately after the email paragraph.
$newPara = $d.Paragraphs.Item($emailParaIndex + 1)
$newRange = $newPara.Range
$newRange.InsertBefore("ShankyShako.GitHub.io")

# Ensure the new paragraph's text run AND its paragraph mark carry the same
# bold / 9pt formatting as the email line above it.
$newPara2 = $d.Paragraphs.Item($emailParaIndex + 1)
$newRange2 = $newPara2.Range
$newRange2.Font.Bold = 1
$newRange2.Font.Size = 9
